# ABT_financials_quarterlyinfo.xlsx
# "completed a few more stock tests, end sunday"
#
# Adds a "sell"-signal mirror of the existing buy-signal analysis
# (columns S/T/U) into new columns X/Y/Z:
#   X{r} = IF(M{r}<0, 1, 0)                     -- a "sell" signal fired
#   Y{r} = IF(AND(M{r}<0, I{r+1}<0), 1, 0)       -- the signal was "correct"
#   X2   = "sell"                               (header-ish label, mirrors J/K cols)
#   X37/Y37/Z37 = "sell times" / "correct" / "percentage"  (summary labels)
#   X38 = SUM(X3:X30)  "sell times"
#   Y38 = SUM(Y3:Y30)  "correct"
#   Z38 = Y38/X38      "percentage"
# Correct/incorrect sell calls are highlighted green/red, matching the
# existing look of the buy-signal columns (S/T/U).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$greenColor = 5296274   # RGB(146,208,80) -- same fill used by correct "buy" calls
$redColor   = 255       # RGB(255,0,0)    -- same fill used by incorrect "buy" calls

# Header-ish label in row 2 (mirrors existing S2/T2/U2 labels)
$ws.Range("X2").Value = "sell"

# Per-row sell signal + correctness flag, rows 3..30
for ($r = 3; $r -le 30; $r++) {
    $next = $r + 1

    $xCell = $ws.Range("X$r")
    $xCell.Formula = "=IF(M$r<0,1,0)"

    $yCell = $ws.Range("Y$r")
    $yCell.Formula = "=IF(AND(M$r<0,I$next<0),1,0)"

    $xVal = $xCell.Value2
    $yVal = $yCell.Value2

    if ($xVal -eq 1) {
        if ($yVal -eq 1) {
            $xCell.Interior.Color = $greenColor
            $yCell.Interior.Color = $greenColor
        } else {
            $xCell.Interior.Color = $redColor
            $yCell.Interior.Color = $redColor
        }
    }
}

# Summary labels, row 37
$ws.Range("X37").Value = "sell times"
$ws.Range("Y37").Value = "correct"
$ws.Range("Z37").Value = "percentage"

# Summary totals, row 38
$ws.Range("X38").Formula = "=SUM(X3:X30)"
$ws.Range("Y38").Formula = "=SUM(Y3:Y30)"
$ws.Range("Z38").Formula = "=Y38/X38"

# View state: zoomed to 70%, selection parked on the new grand-total cell
$excel.ActiveWindow.Zoom = 70
$ws.Range("X38").Select()
